# Add a new "MVerb" (Reverb) effect block to the parameter table.
#
# The sheet lists one row per effect parameter (Project / descriptive name /
# base name / label / enum count / unit / min / max / default / step /
# variable name). Existing data occupies rows 1-77 (row 78 is left blank as
# a separator, matching the blank-row pattern already used between every
# other effect block). We append 8 new rows (79-86) describing the Reverb
# effect's parameters.
#
# NOTE on write order: cell writes below are intentionally sequenced so that
# each *new* distinct text value is written for the first time in the same
# order the source workbook introduced it into the shared-string table
# (Reverb; Density; Bandwidth; Decay; Predelay; Size; Mix; density;
# predelay; roomSize; dryWetMix; damping; bandwidth; Room Size; Damping;
# Late vs Early Mix; Early/Late Mix; earlyLateMix). This keeps the
# regenerated xl/sharedStrings.xml table identical to the target. Cells
# that repeat an already-introduced string (or hold numbers) are filled
# afterwards in straightforward reading order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Phase 1: introduce every new shared string in the required order ---
$ws.Range("A79").Value = "Reverb"

$ws.Range("B80").Value = "Density"
$ws.Range("B81").Value = "Bandwidth"
$ws.Range("B82").Value = "Decay"
$ws.Range("B83").Value = "Predelay"
$ws.Range("B84").Value = "Size"
$ws.Range("B85").Value = "Mix"

$ws.Range("C80").Value = "density"
$ws.Range("C83").Value = "predelay"
$ws.Range("C84").Value = "roomSize"
$ws.Range("C85").Value = "dryWetMix"

$ws.Range("C79").Value = "damping"
$ws.Range("C81").Value = "bandwidth"

$ws.Range("D84").Value = "Room Size"

$ws.Range("B79").Value = "Damping"

$ws.Range("B86").Value = "Late vs Early Mix"
$ws.Range("D86").Value = "Early/Late Mix"
$ws.Range("C86").Value = "earlyLateMix"

# --- Phase 2: fill in the remaining cells for each new row ---

# Row 79: Damping
$ws.Range("A79").Value = "Reverb"
$ws.Range("D79").Value = "Damping"
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 1
$ws.Range("I79").Value = 0.5
$ws.Range("J79").Value = 0.01
$ws.Range("K79").Value = "damping"

# Row 80: Density
$ws.Range("A80").Value = "Reverb"
$ws.Range("D80").Value = "Density"
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 1
$ws.Range("I80").Value = 0.75
$ws.Range("J80").Value = 0.01
$ws.Range("K80").Value = "density"

# Row 81: Bandwidth
$ws.Range("A81").Value = "Reverb"
$ws.Range("D81").Value = "Bandwidth"
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 1
$ws.Range("I81").Value = 1
$ws.Range("J81").Value = 0.01
$ws.Range("K81").Value = "bandwidth"

# Row 82: Decay (reuses existing "decay" base-name string from Flanger/AutoWah)
$ws.Range("A82").Value = "Reverb"
$ws.Range("C82").Value = "decay"
$ws.Range("D82").Value = "Decay"
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 1
$ws.Range("I82").Value = 0.75
$ws.Range("J82").Value = 0.01
$ws.Range("K82").Value = "decay"

# Row 83: Predelay
$ws.Range("A83").Value = "Reverb"
$ws.Range("D83").Value = "Predelay"
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 1
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0.01
$ws.Range("K83").Value = "predelay"

# Row 84: Size (labelled "Room Size")
$ws.Range("A84").Value = "Reverb"
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 1
$ws.Range("I84").Value = 0.8
$ws.Range("J84").Value = 0.01
$ws.Range("K84").Value = "roomSize"

# Row 85: Mix
$ws.Range("A85").Value = "Reverb"
$ws.Range("D85").Value = "Mix"
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 1
$ws.Range("I85").Value = 0.5
$ws.Range("J85").Value = 0.01
$ws.Range("K85").Value = "dryWetMix"

# Row 86: Late vs Early Mix (labelled "Early/Late Mix")
$ws.Range("A86").Value = "Reverb"
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 1
$ws.Range("I86").Value = 0.5
$ws.Range("J86").Value = 0.01
$ws.Range("K86").Value = "earlyLateMix"

# Match the author's final selection (scrolled further down the sheet).
$ws.Range("K90").Select()
